$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0. Drop the two pre-existing hyperlinks before shuffling columns -- the
#    engine doesn't re-anchor Hyperlink ranges when a column is inserted,
#    so we recreate every hyperlink from scratch once the layout settles.
$ws.Hyperlinks.Delete()

# 1. Insert a new column before column A, shifting old A->B, B->C, C->D
$ws.Columns.Item(1).Insert()

# 2. Fill the new column A with its "Label"/"score" tags (rows 2-9)
$ws.Range("A2").Value = "score"
$ws.Range("A3").Value = "score"
$ws.Range("A4").Value = "Label"
$ws.Range("A5").Value = "score"
$ws.Range("A6").Value = "score"
$ws.Range("A7").Value = "score"
$ws.Range("A8").Value = "score"
$ws.Range("A9").Value = "score"

# 3. Fill in the new article name column (C) for the newly-populated rows 4-9
$ws.Range("C4").Value = "Supplier’s Efficiency and Performance  Evaluation using DEA-SVM Approach"
$ws.Range("C5").Value = "A combined neural network and DEA for measuring efficiency of large scale datasets"
$ws.Range("C6").Value = "Farm efficiency estimation using a hybrid approach of machine-learning and data envelopment analysis: Evidence from rural eastern India"
$ws.Range("C7").Value = "Assessing countries’ performances against COVID-19 via WSIDEA and machine learning algorithms"
$ws.Range("C8").Value = "Efficiency analysis for stochastic dynamic facility layout problem using meta-heuristic, data envelopment analysis and machine learning"
$ws.Range("C9").Value = "Using inverse DEA and machine learning algorithms to evaluate and predict suppliers’ performance in the apple supply chain"

# 4. Re-add every hyperlink at its (now correct, post-shift) address: the
#    original two plus the six new ones.
$ws.Hyperlinks.Add($ws.Range("D2"), "https://link.springer.com/article/10.1007/s10479-023-05230-8")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.sciencedirect.com/science/article/pii/S2096232020300469?via%3Dihub")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.jsoftware.us/vol8/jsw0801-04.pdf")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.sciencedirect.com/science/article/pii/S0360835208001113")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.sciencedirect.com/science/article/pii/S0959652620321533", "fig2")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.sciencedirect.com/science/article/pii/S1568494620307304")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://onlinelibrary.wiley.com/doi/full/10.1111/coin.12251")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.sciencedirect.com/science/article/pii/S0925527324000604", "sec3")

# Reuse the workbook's existing "Hipervínculo" style for all hyperlink cells (D2:D9)
$ws.Range("D2:D9").Style = "Hipervínculo"

# 5. Formatting touch-ups: wrap text on the new C4 description, yellow highlight on new A4 tag
$ws.Range("C4").WrapText = $true
$ws.Range("A4").Interior.Color = 65535

# 6. Leave the selection where the author left it
$ws.Range("C10").Select()
